$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to hold a literal string (matches the source file, which
    # stores every data value as inline text) instead of letting Excel
    # auto-convert number-looking strings ("1.002", "243.56", ...) to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.441.56"
$ws.Range("E2").Value = "  +0.48%  "

Set-TextValue $ws.Range("D3") "1.883.45"
$ws.Range("E3").Value = "  +0.54%  "

Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  +0.17%  "

Set-TextValue $ws.Range("D5") "0.7194"
$ws.Range("E5").Value = "  +1.47%  "

Set-TextValue $ws.Range("D6") "243.56"
$ws.Range("E6").Value = "  +0.74%  "

Set-TextValue $ws.Range("D7") "1.002"
$ws.Range("E7").Value = "  +0.16%  "

Set-TextValue $ws.Range("D8") "0.07957"
$ws.Range("E8").Value = "  +2.11%  "

Set-TextValue $ws.Range("D9") "0.3153"
$ws.Range("E9").Value = "  +1.32%  "

Set-TextValue $ws.Range("D10") "24.97"
$ws.Range("E10").Value = "  -0.41%  "

Set-TextValue $ws.Range("D11") "0.08157"
$ws.Range("E11").Value = "  -2.82%  "

Set-TextValue $ws.Range("D12") "1.904.31"
$ws.Range("E12").Value = "  +1.40%  "

Set-TextValue $ws.Range("D15") "0.7111"
$ws.Range("E15").Value = "  -0.86%  "

Set-TextValue $ws.Range("D16") "6.400"
$ws.Range("E16").Value = "  +4.84%  "

Set-TextValue $ws.Range("D17") "0.000008434"
$ws.Range("E17").Value = "  +1.39%  "

Set-TextValue $ws.Range("D18") "29.444.68"

Set-TextValue $ws.Range("D19") "252.47"
$ws.Range("E19").Value = "  +4.98%  "

Set-TextValue $ws.Range("D20") "13.34"
$ws.Range("E20").Value = "  +0.93%  "

Set-TextValue $ws.Range("D21") "2.141.20"
$ws.Range("E21").Value = "  +0.99%  "

$ws.Range("E22").Value = "  +0.13%  "

Set-TextValue $ws.Range("D23") "7.782"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("E24").Value = "  +0.13%  "

Set-TextValue $ws.Range("D25") "0.1589"
$ws.Range("E25").Value = "  +0.00%  "

Set-TextValue $ws.Range("D26") "9.078"
$ws.Range("E26").Value = "  +0.48%  "

Set-TextValue $ws.Range("D27") "162.47"
$ws.Range("E27").Value = "  +0.11%  "

Set-TextValue $ws.Range("D28") "18.96"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("E29").Value = "  +0.25%  "

Set-TextValue $ws.Range("D30") "4.425"
$ws.Range("E30").Value = "  +0.41%  "

Set-TextValue $ws.Range("D31") "4.287"
$ws.Range("E31").Value = "  -0.66%  "

Set-TextValue $ws.Range("D32") "1.225"
$ws.Range("E32").Value = "  -3.34%  "

Set-TextValue $ws.Range("D33") "0.05328"
$ws.Range("E33").Value = "  -0.86%  "

Set-TextValue $ws.Range("D34") "1.947"
$ws.Range("E34").Value = "  +0.43%  "

Set-TextValue $ws.Range("D35") "0.7557"
$ws.Range("E35").Value = "  +0.67%  "

Set-TextValue $ws.Range("D36") "1.181"
$ws.Range("E36").Value = "  +0.47%  "

Set-TextValue $ws.Range("D37") "2.702"
$ws.Range("E37").Value = "  +0.73%  "

Set-TextValue $ws.Range("D38") "0.01885"
$ws.Range("E38").Value = "  +0.61%  "

Set-TextValue $ws.Range("D39") "1.275.11"
$ws.Range("E39").Value = "  +3.05%  "

Set-TextValue $ws.Range("D40") "2.770"
$ws.Range("E40").Value = "  +1.44%  "

Set-TextValue $ws.Range("D41") "6.479"
$ws.Range("E41").Value = "  -0.60%  "

Set-TextValue $ws.Range("D42") "113.17"
$ws.Range("E42").Value = "  +3.77%  "

Set-TextValue $ws.Range("D43") "74.67"
$ws.Range("E43").Value = "  +3.13%  "

Set-TextValue $ws.Range("D44") "0.9085"
$ws.Range("E44").Value = "  +1.91%  "

Set-TextValue $ws.Range("D45") "0.00000000131"
$ws.Range("E45").Value = "  +4.40%  "

$ws.Range("E46").Value = "  +0.14%  "

Set-TextValue $ws.Range("D47") "2.033.63"
$ws.Range("E47").Value = "  +0.66%  "

Set-TextValue $ws.Range("D48") "1.805"
$ws.Range("E48").Value = "  +0.76%  "

Set-TextValue $ws.Range("D49") "0.5198"
$ws.Range("E49").Value = "  +0.01%  "

Set-TextValue $ws.Range("D50") "9.523"
$ws.Range("E50").Value = "  +1.08%  "

Set-TextValue $ws.Range("D51") "0.4372"
$ws.Range("E51").Value = "  +0.76%  "

# Rows 13 and 14 swapped places (Litecoin <-> Polkadot) with updated figures
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "5.246"
$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D14") "94.76"
$ws.Range("E14").Value = "  +3.99%  "

